# Add the three new schedule rows (第85-87期) below the existing data,
# widen columns A:B so the new/longer date text isn't truncated, and
# leave the selection on the new last row of column C (mirrors what a
# user does after typing the last entry and pressing Enter).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New rows -------------------------------------------------------
$ws.Range("A35").Value = "2025/11/21"
$ws.Range("B35").Value = "2026/1/16"
$ws.Range("C35").Value = "第85期 第七代坐騎"

$ws.Range("A36").Value = "2025/11/28"
$ws.Range("B36").Value = "2026/1/23"
$ws.Range("C36").Value = "第86期 第四代寵物"

$ws.Range("A37").Value = "2025/11/35"
$ws.Range("B37").Value = "2026/1/30"
$ws.Range("C37").Value = "第87期 十轉技能"

# --- Widen the date columns so the longer strings fit ---------------
$ws.Range("A:B").ColumnWidth = 14.69921875

# --- Scroll / selection, matching where the author left off ---------
$ws.Range("C38").Select() | Out-Null
